# "micole finish dorm management module"
# Replace the two warden records (warden 2 / warden 3) with fresh ones
# (warden 4 / warden 5), give each a working mailto hyperlink, and move
# the selection cursor.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Names (column A) -----------------------------------------------------
$ws.Range("A2").Value = "warden 4"
$ws.Range("A3").Value = "warden 5"

# --- Emails (column B) -----------------------------------------------------
$ws.Range("B2").Value = "warden4@gmail.com"
$ws.Range("B3").Value = "warden5@gmail.com"

# --- Phone numbers (column C) ----------------------------------------------
$ws.Range("C2").Value = 109654782
$ws.Range("C3").Value = 123698699

# Drop the old hyperlink on B2 (it still points at the retired
# warden3@gmail.com address) so we can re-create it against the new email.
[void]$ws.Range("B2").Hyperlinks.Delete()

[void]$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:warden4@gmail.com")
$ws.Range("B2").Style = "Hyperlink"

# B3 never had a hyperlink before - add one now.
[void]$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:warden5@gmail.com")
$ws.Range("B3").Style = "Hyperlink"

# Leave the selection where the author left it when they saved.
[void]$ws.Range("D12").Select()
